# Weekly update: insert a new record at the top of the data range (row 116),
# pushing the existing rows 116-225 down to 117-226.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 116; this shifts rows 116..225 down to 117..226
# and Excel extends the sheet dimension to A1:R226 automatically.
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(116, 1).Value2  = 3
$ws.Cells.Item(116, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(116, 3).Value2  = "Coquimbo"
$ws.Cells.Item(116, 4).Value2  = 44484
$ws.Cells.Item(116, 5).Value2  = 5
$ws.Cells.Item(116, 6).Value2  = 100112031
$ws.Cells.Item(116, 7).Value2  = "Poroto verde"
$ws.Cells.Item(116, 8).Value2  = "Magnum"
$ws.Cells.Item(116, 9).Value2  = "Primera"
$ws.Cells.Item(116, 10).Value2 = 38
$ws.Cells.Item(116, 11).Value2 = 40000
$ws.Cells.Item(116, 12).Value2 = 40000
$ws.Cells.Item(116, 13).Value2 = 40000
$ws.Cells.Item(116, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(116, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(116, 16).Value2 = 1600
$ws.Cells.Item(116, 17).Value2 = 25
$ws.Cells.Item(116, 18).Value2 = "Hortaliza"
